# RepartoTG2 - "Correcciones, añadido datos excel"
# Fill in the missing "Alberto" assignment for 6.2) Situación 2, and
# move the active selection the way the author left it (E16).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# D14 ("6.2) Situación 2") was left blank - the author assigned it to Alberto.
$ws.Range("D14").Value = "Alberto"

# Match the resulting selection left by the editing session.
$ws.Range("E16").Select()
